$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "64.226.55"
$ws.Range("E2").Value = "  -0.98%  "

$ws.Range("D3").Value = "3.512.49"
$ws.Range("E3").Value = "  -0.31%  "

Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.04%  "

Set-TextValue $ws.Range("D5") "586.08"
$ws.Range("E5").Value = "  -1.78%  "

Set-TextValue $ws.Range("D6") "133.07"
$ws.Range("E6").Value = "  -0.85%  "

$ws.Range("D7").Value = "3.510.85"
$ws.Range("E7").Value = "  -0.33%  "

$ws.Range("E8").Value = "  +0.06%  "

Set-TextValue $ws.Range("D9") "0.489"
$ws.Range("E9").Value = "  -1.28%  "

$ws.Range("E10").Value = "  +0.30%  "

Set-TextValue $ws.Range("D11") "7.16"
$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("E12").Value = "  -0.14%  "

$ws.Range("D13").Value = "4.114.53"
$ws.Range("E13").Value = "  -0.04%  "

Set-TextValue $ws.Range("D14") "27.78"
$ws.Range("E14").Value = "  +1.62%  "

$ws.Range("E15").Value = "  +1.73%  "

$ws.Range("E16").Value = "  -0.85%  "

$ws.Range("D17").Value = "3.526.74"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").Value = "64.250.27"
$ws.Range("E18").Value = "  -0.98%  "

Set-TextValue $ws.Range("D19") "10.12"
$ws.Range("E19").Value = "  +2.02%  "

Set-TextValue $ws.Range("D20") "14.44"
$ws.Range("E20").Value = "  -0.08%  "

Set-TextValue $ws.Range("D21") "5.69"
$ws.Range("E21").Value = "  -0.20%  "

Set-TextValue $ws.Range("D22") "386.61"
$ws.Range("E22").Value = "  -1.36%  "

Set-TextValue $ws.Range("D23") "0.581"
$ws.Range("E23").Value = "  +0.74%  "

$ws.Range("D24").Value = "3.653.23"
$ws.Range("E24").Value = "  -0.27%  "

Set-TextValue $ws.Range("D25") "73.28"
$ws.Range("E25").Value = "  -1.13%  "

$ws.Range("E26").Value = "  +0.02%  "

Set-TextValue $ws.Range("D27") "0.0000115"
$ws.Range("E27").Value = "  +1.71%  "

$ws.Range("E28").Value = "  -2.11%  "

Set-TextValue $ws.Range("D29") "7.64"
$ws.Range("E29").Value = "  -2.05%  "

$ws.Range("E30").Value = "  -0.24%  "

$ws.Range("E31").Value = "  -0.89%  "

Set-TextValue $ws.Range("D32") "8.33"
$ws.Range("E32").Value = "  -0.90%  "

$ws.Range("D33").Value = "3.519.99"
$ws.Range("E33").Value = "  -0.15%  "

Set-TextValue $ws.Range("D35") "23.81"
$ws.Range("E35").Value = "  -1.28%  "

$ws.Range("E36").Value = "  +0.13%  "

$ws.Range("E37").Value = "  +3.18%  "

Set-TextValue $ws.Range("D38") "1.59"
$ws.Range("E38").Value = "  +0.31%  "

$ws.Range("E39").Value = "  +1.91%  "

Set-TextValue $ws.Range("D40") "163.70"
$ws.Range("E40").Value = "  -3.02%  "

Set-TextValue $ws.Range("D41") "0.0805"
$ws.Range("E41").Value = "  -1.94%  "

Set-TextValue $ws.Range("D46") "41.80"
$ws.Range("E46").Value = "  -2.03%  "

$ws.Range("E47").Value = "  -0.25%  "

Set-TextValue $ws.Range("D48") "1.65"
$ws.Range("E48").Value = "  +0.16%  "

Set-TextValue $ws.Range("D49") "6.88"
$ws.Range("E49").Value = "  -0.67%  "

$ws.Range("D50").Value = "2.431.79"
$ws.Range("E50").Value = "  +1.57%  "

# Row re-mappings (entity swaps / replacements)
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D42") "26.66"
$ws.Range("E42").Value = "  +5.19%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D43") "0.814"
$ws.Range("E43").Value = "  -0.80%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D44") "1.00"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D45") "1.23"
$ws.Range("E45").Value = "  -0.89%  "

$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D51") "0.900"
$ws.Range("E51").Value = "  -0.16%  "
